# Medical Cost Descriptions.xlsx - edit script
#
# Commit message: "Created our java exercise for our new protocol, MOSTLY
# finished our add-in sheet w/ descriptions. Still need a calibrated costs
# burden (out of pocket and insurance costs)."
#
# Content change: the "Complementary" row's Basic-plan description (C12 on
# Sheet1) is reworded - "chiropractic (for back, neck or bone problems)" is
# replaced with "spiritual healing or Reiki".
#
# Cosmetic change: the sheet's scroll/selection state moved from
# topLeftCell=B1 / selection F3 to topLeftCell=A4 / selection D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$newText = "1 point: Covers " + $openQuote + "alternative" + $closeQuote + " services including acupuncture (for pain), spiritual healing or Reiki, and therapeutic massage."

$ws.Range("C12").Value = $newText

# Reflect the new scroll position / selection recorded in the saved view.
$ws.Activate()
$ws.Range("D11").Select()
